# Generate Report for Handoff
# Two new files are now "Ready for handoff":
#   0906733a-99ef-4808-b61a-c26b0135c58a
#   1c193648-01af-4ad9-b7a5-0685e7c86164
# They are inserted ahead of the pre-existing 551ef9d2-... row (which is
# pushed down), on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# Row 3 now describes the first new file (was 551ef9d2 before the edit)
$wsOv.Range("A3").Value = "0906733a-99ef-4808-b61a-c26b0135c58a.md"
$wsOv.Range("B3").Value = "Ready for handoff"
$wsOv.Range("C3").Value = "Ready for handoff"
$wsOv.Range("D3").Value = "2016-03-24 14:49:56"

# Row 4 (new): second new file
$wsOv.Range("A4").Value = "1c193648-01af-4ad9-b7a5-0685e7c86164.md"
$wsOv.Range("B4").Value = "Ready for handoff"
$wsOv.Range("C4").Value = "Ready for handoff"
$wsOv.Range("D4").Value = "2016-03-24 14:49:56"

# Row 5 (new): the original 551ef9d2 row, now pushed down
$wsOv.Range("A5").Value = "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md"
$wsOv.Range("B5").Value = "Ready for handoff"
$wsOv.Range("C5").Value = "Ready for handoff"
$wsOv.Range("D5").Value = "2016-03-24 14:47:32"

# Rebuild hyperlinks top to bottom (A2..A5) so ids line up and no range is
# registered twice.
$wsOv.Hyperlinks.Delete()
$wsOv.Hyperlinks.Add($wsOv.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a387dd077c39973a1e1ecc886b6ad9d2393b080b/e2e/f88f1c75-83c7-44b2-be9d-341104fcb25d.md", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.md") | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9bd212f757e9da0c4a364a8756710465f1e9d0d1/e2e/0906733a-99ef-4808-b61a-c26b0135c58a.md", "", "", "0906733a-99ef-4808-b61a-c26b0135c58a.md") | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b20529df32e237380de79c28bdca4b900b76fdd1/e2e/1c193648-01af-4ad9-b7a5-0685e7c86164.md", "", "", "1c193648-01af-4ad9-b7a5-0685e7c86164.md") | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/3407373e5ef3d9cd09a1eba61467cdb0041a0c9e/e2e/551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md", "", "", "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status |
#   Latest Handoff File | Latest Handoff Datetime | Latest Target File |
#   Latest Handback File | Latest Handback DateTime | Reference Tokens |
#   Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = "0906733a-99ef-4808-b61a-c26b0135c58a.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "0906733a-99ef-4808-b61a-c26b0135c58a.813da62001f07b367304767492ec430e95058069.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-24 14:49:52"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("J3").Value = "Include"

$wsZh.Range("A4").Value = "1c193648-01af-4ad9-b7a5-0685e7c86164.md"
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "1c193648-01af-4ad9-b7a5-0685e7c86164.852bfc55ab6c2ad83c68fa01d85124ea0c4872df.zh-cn.xlf"
$wsZh.Range("E4").Value = "2016-03-24 14:49:52"
$wsZh.Range("H4").Value = "0001-01-01 00:00:00"
$wsZh.Range("J4").Value = "Include"

$wsZh.Range("A5").Value = "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md"
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.460871851b3c3f69f4cdb5f568904f746a515d1c.zh-cn.xlf"
$wsZh.Range("E5").Value = "2016-03-24 14:47:27"
$wsZh.Range("H5").Value = "0001-01-01 00:00:00"
$wsZh.Range("J5").Value = "Include"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a387dd077c39973a1e1ecc886b6ad9d2393b080b/e2e/f88f1c75-83c7-44b2-be9d-341104fcb25d.md", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a3384ba76bdf6b091b306fb3c250dab4e483a9a5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.zh-cn.xlf", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e28f19bdc274a5066b0708a5ba7bc25277a42691/e2e/f88f1c75-83c7-44b2-be9d-341104fcb25d.md", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/54493acffd0a7c0fa1c9fbb0a24d0bc039f8258c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.zh-cn.xlf", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9bd212f757e9da0c4a364a8756710465f1e9d0d1/e2e/0906733a-99ef-4808-b61a-c26b0135c58a.md", "", "", "0906733a-99ef-4808-b61a-c26b0135c58a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/edc923877bb0da5fbd50e6299b3f6dc0c3fd4ec3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0906733a-99ef-4808-b61a-c26b0135c58a.813da62001f07b367304767492ec430e95058069.zh-cn.xlf", "", "", "0906733a-99ef-4808-b61a-c26b0135c58a.813da62001f07b367304767492ec430e95058069.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b20529df32e237380de79c28bdca4b900b76fdd1/e2e/1c193648-01af-4ad9-b7a5-0685e7c86164.md", "", "", "1c193648-01af-4ad9-b7a5-0685e7c86164.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3ea13a9866632c77ad84bb4b05433ca2f9268586/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1c193648-01af-4ad9-b7a5-0685e7c86164.852bfc55ab6c2ad83c68fa01d85124ea0c4872df.zh-cn.xlf", "", "", "1c193648-01af-4ad9-b7a5-0685e7c86164.852bfc55ab6c2ad83c68fa01d85124ea0c4872df.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/3407373e5ef3d9cd09a1eba61467cdb0041a0c9e/e2e/551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md", "", "", "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e148cd17d31df31b4a8823697e7085eaa0ed62ef/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/551ef9d2-4f53-4ab1-9a3a-b25095c949aa.460871851b3c3f69f4cdb5f568904f746a515d1c.zh-cn.xlf", "", "", "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.460871851b3c3f69f4cdb5f568904f746a515d1c.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de": same layout as "zh-cn"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = "0906733a-99ef-4808-b61a-c26b0135c58a.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "0906733a-99ef-4808-b61a-c26b0135c58a.813da62001f07b367304767492ec430e95058069.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-24 14:49:56"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("J3").Value = "Include"

$wsDe.Range("A4").Value = "1c193648-01af-4ad9-b7a5-0685e7c86164.md"
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "1c193648-01af-4ad9-b7a5-0685e7c86164.852bfc55ab6c2ad83c68fa01d85124ea0c4872df.de-de.xlf"
$wsDe.Range("E4").Value = "2016-03-24 14:49:56"
$wsDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDe.Range("J4").Value = "Include"

$wsDe.Range("A5").Value = "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md"
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.460871851b3c3f69f4cdb5f568904f746a515d1c.de-de.xlf"
$wsDe.Range("E5").Value = "2016-03-24 14:47:32"
$wsDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDe.Range("J5").Value = "Include"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a387dd077c39973a1e1ecc886b6ad9d2393b080b/e2e/f88f1c75-83c7-44b2-be9d-341104fcb25d.md", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f5e5c43f3582d1630d5eb7ab1143cd42aeca3b62/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.de-de.xlf", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/068eed81d08c0e6118b130dcc180bd781c7847c0/e2e/f88f1c75-83c7-44b2-be9d-341104fcb25d.md", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/aee2c4d239dab3a9c68de2dc5fcca3a964433e36/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.de-de.xlf", "", "", "f88f1c75-83c7-44b2-be9d-341104fcb25d.8a58f7b110dc4945d83aecf21203bf2d9247b49a.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9bd212f757e9da0c4a364a8756710465f1e9d0d1/e2e/0906733a-99ef-4808-b61a-c26b0135c58a.md", "", "", "0906733a-99ef-4808-b61a-c26b0135c58a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ef80f3b51dc6f7162c09ed0a5d76b981788a1c2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0906733a-99ef-4808-b61a-c26b0135c58a.813da62001f07b367304767492ec430e95058069.de-de.xlf", "", "", "0906733a-99ef-4808-b61a-c26b0135c58a.813da62001f07b367304767492ec430e95058069.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b20529df32e237380de79c28bdca4b900b76fdd1/e2e/1c193648-01af-4ad9-b7a5-0685e7c86164.md", "", "", "1c193648-01af-4ad9-b7a5-0685e7c86164.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/56d140a92be700b29b2e6b89a76021098cacf6a6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1c193648-01af-4ad9-b7a5-0685e7c86164.852bfc55ab6c2ad83c68fa01d85124ea0c4872df.de-de.xlf", "", "", "1c193648-01af-4ad9-b7a5-0685e7c86164.852bfc55ab6c2ad83c68fa01d85124ea0c4872df.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/3407373e5ef3d9cd09a1eba61467cdb0041a0c9e/e2e/551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md", "", "", "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1ddce9b2df8afb83fc283938c526c02aba70f64c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/551ef9d2-4f53-4ab1-9a3a-b25095c949aa.460871851b3c3f69f4cdb5f568904f746a515d1c.de-de.xlf", "", "", "551ef9d2-4f53-4ab1-9a3a-b25095c949aa.460871851b3c3f69f4cdb5f568904f746a515d1c.de-de.xlf") | Out-Null
